# Add date support for jobs: insert a "Name" column right after "Id",
# shift "Client"/"Type" over, and append a new "Date" column with a
# fixed timestamp for every job row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateValue = "Feb 6, 2022 (05:00:11 EST)"

# Header row
$ws.Cells.Item(1, 1).Value = "Id"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "Client"
$ws.Cells.Item(1, 4).Value = "Type"
$ws.Cells.Item(1, 5).Value = "File"
$ws.Cells.Item(1, 6).Value = "Date"

# Data rows: Id, Name, Client, Type, File, Date
$rows = @(
    @(1, "Supplier",     "DRX", "I", "I_Supplier.xml"),
    @(2, "Plant",        "DRX", "I", "I_Plant.xml"),
    @(3, "Solicitation", "DRX", "I", "I_Solicitation.xml"),
    @(4, "BOM",           "DRX", "I", "I_BOM.xml"),
    @(5, "RequestFile",  "DRX", "E", "E_RequestFile.xml"),
    @(6, "BOM",           "DRX", "E", "E_BOM.xml"),
    @(7, "Supplier",     "GYU", "I", "I_Supplier.xml"),
    @(8, "Plant",        "GYU", "I", "I_Plant.xml"),
    @(9, "Solicitation", "GYU", "I", "I_Solicitation.xml"),
    @(10, "BOM",          "GYU", "I", "I_BOM.xml"),
    @(11, "RequestFile", "GYU", "E", "E_RequestFile.xml"),
    @(12, "BOM",          "GYU", "E", "E_BOM.xml")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = [double]$row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $dateValue
    $r = $r + 1
}
